$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 69
$ws.Range("F3").Value = 3321
$ws.Range("F5").Value = 1366
$ws.Range("F7").Value = 3826
$ws.Range("F10").Value = 55
$ws.Range("F11").Value = 8551
$ws.Range("F12").Value = 8551
$ws.Range("F13").Value = 475
$ws.Range("F16").Value = 103
$ws.Range("F17").Value = 325
$ws.Range("F19").Value = 81
$ws.Range("F20").Value = 8
$ws.Range("F22").Value = 10797
$ws.Range("F23").Value = 10797
$ws.Range("F24").Value = 288
$ws.Range("F26").Value = 27
$ws.Range("F30").Value = 158
$ws.Range("F31").Value = 136
$ws.Range("F32").Value = 2668
$ws.Range("F34").Value = 94
$ws.Range("F38").Value = 2118
$ws.Range("F40").Value = 4074
$ws.Range("F41").Value = 2161
$ws.Range("F43").Value = 2577
$ws.Range("F44").Value = 3023
$ws.Range("F45").Value = 1237
$ws.Range("F46").Value = 168
$ws.Range("F48").Value = 337
$ws.Range("F49").Value = 305
$ws.Range("F51").Value = 118

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 13
$ws.Range("F6").Value = 189
$ws.Range("F7").Value = 43
$ws.Range("F9").Value = 4
$ws.Range("F21").Value = 29

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 20

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 69
$ws.Range("F5").Value = 3321
$ws.Range("F6").Value = 189
$ws.Range("F8").Value = 1366
$ws.Range("F10").Value = 3826
$ws.Range("F12").Value = 43
$ws.Range("F16").Value = 8551
$ws.Range("F17").Value = 475
$ws.Range("F20").Value = 103
$ws.Range("F21").Value = 325
$ws.Range("F23").Value = 81
$ws.Range("F24").Value = 8
$ws.Range("F25").Value = 10797
$ws.Range("F26").Value = 288
$ws.Range("F27").Value = 27
$ws.Range("F28").Value = 20
$ws.Range("F33").Value = 158
$ws.Range("F34").Value = 136
$ws.Range("F35").Value = 2668
$ws.Range("F37").Value = 94
$ws.Range("F41").Value = 2118
$ws.Range("F44").Value = 2161
$ws.Range("F45").Value = 3023
$ws.Range("F47").Value = 1237
$ws.Range("F48").Value = 337
$ws.Range("F49").Value = 305
$ws.Range("F51").Value = 118

$wb.Save()
